$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.420.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.63%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.616.08'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.49%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.500'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.245'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0608'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.843.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.622.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.82%  '

$ws.Range("E14").Value = '  +0.48%  '

$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '237.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +10.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.432.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0725'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.33%  '

$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.96%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.52%  '

$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("E28").Value = '  +0.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.47%  '

$ws.Range("E30").Value = '  +0.34%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.526.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.53%  '

$ws.Range("E33").Value = '  +1.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.56%  '

$ws.Range("E35").Value = '  +4.46%  '

$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.565'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0166'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.48%  '

$ws.Range("E39").Value = '  +0.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.42%  '

$ws.Range("E41").Value = '  -0.05%  '

$ws.Range("E42").Value = '  +1.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.755.14'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.764'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.908'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.22%  '

$ws.Range("E48").Value = '  +1.96%  '

$ws.Range("E49").Value = '  +0.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0960'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.17%  '
